# Update recomputed TPM-based NATMI edge metrics for Egf-Egfr (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1463976666666667
$ws.Range("H2").Value = 0.439193
$ws.Range("I2").Value = 0.1157910139257259
$ws.Range("J2").Value = 0.115791013925726
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 0.2491585320307778
$ws.Range("R2").Value = 2.242426788277
$ws.Range("S2").Value = 0.002438122090465728
$ws.Range("T2").Value = 0.002438122090465728

# Row 3
$ws.Range("G3").Value = 0.1463976666666667
$ws.Range("H3").Value = 0.439193
$ws.Range("I3").Value = 0.1157910139257259
$ws.Range("J3").Value = 0.115791013925726
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 9.150431998257002
$ws.Range("R3").Value = 82.353887984313
$ws.Range("S3").Value = 0.0895408646471676
$ws.Range("T3").Value = 0.08954086464716761

# Row 4
$ws.Range("G4").Value = 0.1463976666666667
$ws.Range("H4").Value = 0.439193
$ws.Range("I4").Value = 0.1157910139257259
$ws.Range("J4").Value = 0.115791013925726
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 0.06351740923900001
$ws.Range("R4").Value = 0.571656683151
$ws.Range("S4").Value = 0.000621544834658233
$ws.Range("T4").Value = 0.0006215448346582332

# Row 5
$ws.Range("G5").Value = 0.1463976666666667
$ws.Range("H5").Value = 0.439193
$ws.Range("I5").Value = 0.1157910139257259
$ws.Range("J5").Value = 0.115791013925726
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 2.338319748704
$ws.Range("R5").Value = 21.044877738336
$ws.Range("S5").Value = 0.0228814521719178
$ws.Range("T5").Value = 0.02288145217191781

# Row 6
$ws.Range("G6").Value = 0.1463976666666667
$ws.Range("H6").Value = 0.439193
$ws.Range("I6").Value = 0.1157910139257259
$ws.Range("J6").Value = 0.115791013925726
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 0.03158066065722222
$ws.Range("R6").Value = 0.284225945915
$ws.Range("S6").Value = 0.0003090301815165782
$ws.Range("T6").Value = 0.0003090301815165782

# Row 7
$ws.Range("I7").Value = 0.4041732358198567
$ws.Range("J7").Value = 0.4041732358198568
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 0.8696979731743335
$ws.Range("R7").Value = 7.827281758569001
$ws.Range("S7").Value = 0.008510364157096904
$ws.Range("T7").Value = 0.008510364157096906

# Row 8
$ws.Range("I8").Value = 0.4041732358198567
$ws.Range("J8").Value = 0.4041732358198568
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.312546023871659
$ws.Range("T8").Value = 0.3125460238716591

# Row 9
$ws.Range("I9").Value = 0.4041732358198567
$ws.Range("J9").Value = 0.4041732358198568
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 0.2217100960830001
$ws.Range("R9").Value = 1.995390864747
$ws.Range("S9").Value = 0.00216952748329914
$ws.Range("T9").Value = 0.00216952748329914

# Row 10
$ws.Range("I10").Value = 0.4041732358198567
$ws.Range("J10").Value = 0.4041732358198568
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 8.162000030688002
$ws.Range("R10").Value = 73.45800027619201
$ws.Range("S10").Value = 0.0798686379109995
$ws.Range("T10").Value = 0.07986863791099953

# Row 11
$ws.Range("I11").Value = 0.4041732358198567
$ws.Range("J11").Value = 0.4041732358198568
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 0.1102335783616667
$ws.Range("R11").Value = 0.9921022052550001
$ws.Range("S11").Value = 0.001078682396802149
$ws.Range("T11").Value = 0.001078682396802149

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1328766666666667
$ws.Range("H12").Value = 0.39863
$ws.Range("I12").Value = 0.1050967840589721
$ws.Range("J12").Value = 0.1050967840589721
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 0.2261467410077778
$ws.Range("R12").Value = 2.03532066907
$ws.Range("S12").Value = 0.002212941938788535
$ws.Range("T12").Value = 0.002212941938788535

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1328766666666667
$ws.Range("H13").Value = 0.39863
$ws.Range("I13").Value = 0.1050967840589721
$ws.Range("J13").Value = 0.1050967840589721
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 8.30531613087
$ws.Range("R13").Value = 74.74784517783
$ws.Range("S13").Value = 0.08127104683886222
$ws.Range("T13").Value = 0.08127104683886223

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1328766666666667
$ws.Range("H14").Value = 0.39863
$ws.Range("I14").Value = 0.1050967840589721
$ws.Range("J14").Value = 0.1050967840589721
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 0.05765106649000001
$ws.Range("R14").Value = 0.5188595984100001
$ws.Range("S14").Value = 0.000564140178554329
$ws.Range("T14").Value = 0.0005641401785543291

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1328766666666667
$ws.Range("H15").Value = 0.39863
$ws.Range("I15").Value = 0.1050967840589721
$ws.Range("J15").Value = 0.1050967840589721
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 2.12235714464
$ws.Range("R15").Value = 19.10121430176
$ws.Range("S15").Value = 0.020768166339836
$ws.Range("T15").Value = 0.020768166339836

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1328766666666667
$ws.Range("H16").Value = 0.39863
$ws.Range("I16").Value = 0.1050967840589721
$ws.Range("J16").Value = 0.1050967840589721
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 0.02866393307222222
$ws.Range("R16").Value = 0.25797539765
$ws.Range("S16").Value = 0.0002804887629309974
$ws.Range("T16").Value = 0.0002804887629309974

# Row 17
$ws.Range("G17").Value = 0.180116
$ws.Range("H17").Value = 0.540348
$ws.Range("I17").Value = 0.1424600182442301
$ws.Range("J17").Value = 0.1424600182442301
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 0.3065447638413333
$ws.Range("R17").Value = 2.758902874572
$ws.Range("S17").Value = 0.002999670749167166
$ws.Range("T17").Value = 0.002999670749167166

# Row 18
$ws.Range("G18").Value = 0.180116
$ws.Range("H18").Value = 0.540348
$ws.Range("I18").Value = 0.1424600182442301
$ws.Range("J18").Value = 0.1424600182442301
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 11.257960917852
$ws.Range("R18").Value = 101.321648260668
$ws.Range("S18").Value = 0.1101639305051941
$ws.Range("T18").Value = 0.1101639305051941

# Row 19
$ws.Range("G19").Value = 0.180116
$ws.Range("H19").Value = 0.540348
$ws.Range("I19").Value = 0.1424600182442301
$ws.Range("J19").Value = 0.1424600182442301
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 0.07814674880400002
$ws.Range("R19").Value = 0.7033207392360001
$ws.Range("S19").Value = 0.0007646991375497946
$ws.Range("T19").Value = 0.0007646991375497947

# Row 20
$ws.Range("G20").Value = 0.180116
$ws.Range("H20").Value = 0.540348
$ws.Range("I20").Value = 0.1424600182442301
$ws.Range("J20").Value = 0.1424600182442301
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 2.876881916544
$ws.Range("R20").Value = 25.891937248896
$ws.Range("S20").Value = 0.02815151179137974
$ws.Range("T20").Value = 0.02815151179137974

# Row 21
$ws.Range("G21").Value = 0.180116
$ws.Range("H21").Value = 0.540348
$ws.Range("I21").Value = 0.1424600182442301
$ws.Range("J21").Value = 0.1424600182442301
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 0.03885432332666668
$ws.Range("R21").Value = 0.3496889099400001
$ws.Range("S21").Value = 0.0003802060609393137
$ws.Range("T21").Value = 0.0003802060609393137

# Row 22
$ws.Range("G22").Value = 0.2939293333333333
$ws.Range("H22").Value = 0.881788
$ws.Range("I22").Value = 0.2324789479512151
$ws.Range("J22").Value = 0.2324789479512152
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 0.5002470523035555
$ws.Range("R22").Value = 4.502223470732
$ws.Range("S22").Value = 0.004895129935831384
$ws.Range("T22").Value = 0.004895129935831385

# Row 23
$ws.Range("G23").Value = 0.2939293333333333
$ws.Range("H23").Value = 0.881788
$ws.Range("I23").Value = 0.2324789479512151
$ws.Range("J23").Value = 0.2324789479512152
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 18.371743472412
$ws.Range("R23").Value = 165.345691251708
$ws.Range("S23").Value = 0.1797753150790122
$ws.Range("T23").Value = 0.1797753150790123

# Row 24
$ws.Range("G24").Value = 0.2939293333333333
$ws.Range("H24").Value = 0.881788
$ws.Range("I24").Value = 0.2324789479512151
$ws.Range("J24").Value = 0.2324789479512152
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 0.127526825924
$ws.Range("R24").Value = 1.147741433316
$ws.Range("S24").Value = 0.001247904171204035
$ws.Range("T24").Value = 0.001247904171204036

# Row 25
$ws.Range("G25").Value = 0.2939293333333333
$ws.Range("H25").Value = 0.881788
$ws.Range("I25").Value = 0.2324789479512151
$ws.Range("J25").Value = 0.2324789479512152
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 4.694752180864
$ws.Range("R25").Value = 42.252769627776
$ws.Range("S25").Value = 0.04594014464659285
$ws.Range("T25").Value = 0.04594014464659286

# Row 26
$ws.Range("G26").Value = 0.2939293333333333
$ws.Range("H26").Value = 0.881788
$ws.Range("I26").Value = 0.2324789479512151
$ws.Range("J26").Value = 0.2324789479512152
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 0.06340594590444444
$ws.Range("R26").Value = 0.57065351314
$ws.Range("S26").Value = 0.000620454118574614
$ws.Range("T26").Value = 0.0006204541185746141

